$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.694.35"
$ws.Range("E2").Value = "'  +1.82%  "
$ws.Range("D3").Value = "'1.897.60"
$ws.Range("E3").Value = "'  +2.78%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("D5").Value = "'239.01"
$ws.Range("E5").Value = "'  +1.28%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = "'  +0.05%  "
$ws.Range("D7").Value = "'0.4822"
$ws.Range("E7").Value = "'  +0.94%  "
$ws.Range("D8").Value = "'0.2847"
$ws.Range("E8").Value = "'  +1.58%  "
$ws.Range("D9").Value = "'0.06550"
$ws.Range("E9").Value = "'  +1.26%  "
$ws.Range("D10").Value = "'2.022.30"
$ws.Range("E10").Value = "'  +9.46%  "
$ws.Range("D11").Value = "'0.07454"
$ws.Range("E11").Value = "'  +1.91%  "
$ws.Range("D12").Value = "'16.69"
$ws.Range("E12").Value = "'  +2.71%  "
$ws.Range("D13").Value = "'5.103"
$ws.Range("E13").Value = "'  -0.03%  "
$ws.Range("D14").Value = "'88.08"
$ws.Range("E14").Value = "'  +1.15%  "
$ws.Range("D15").Value = "'0.6666"
$ws.Range("E15").Value = "'  +3.42%  "
$ws.Range("D16").Value = "'30.675.78"
$ws.Range("E16").Value = "'  +1.97%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D17").Value = "'13.30"
$ws.Range("E17").Value = "'  +0.82%  "
$ws.Range("B18").Value = "Dai"
$ws.Range("C18").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D18").Value = "'1.000"
$ws.Range("E18").Value = "'  -0.10%  "
$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").Value = "'2.206.74"
$ws.Range("E19").Value = "'  +5.18%  "
$ws.Range("D20").Value = "'0.000007611"
$ws.Range("E20").Value = "'  -0.05%  "
$ws.Range("D21").Value = "'231.39"
$ws.Range("E21").Value = "'  +3.33%  "
$ws.Range("B22").Value = "BinanceUSD"
$ws.Range("C22").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "'  +0.04%  "
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "'5.284"
$ws.Range("E23").Value = "'  +0.06%  "
$ws.Range("D24").Value = "'6.231"
$ws.Range("E24").Value = "'  +2.68%  "
$ws.Range("D25").Value = "'169.86"
$ws.Range("E25").Value = "'  +3.99%  "
$ws.Range("D26").Value = "'9.336"
$ws.Range("E26").Value = "'  +1.40%  "
$ws.Range("D27").Value = "'18.75"
$ws.Range("E27").Value = "'  +1.51%  "
$ws.Range("D28").Value = "'1.967"
$ws.Range("E28").Value = "'  +2.84%  "
$ws.Range("D29").Value = "'1.404"
$ws.Range("E29").Value = "'  -1.70%  "
$ws.Range("D30").Value = "'0.1017"
$ws.Range("E30").Value = "'  +10.68%  "
$ws.Range("D31").Value = "'4.358"
$ws.Range("E31").Value = "'  +2.95%  "
$ws.Range("D32").Value = "'4.027"
$ws.Range("E32").Value = "'  +1.95%  "
$ws.Range("D33").Value = "'0.05122"
$ws.Range("E33").Value = "'  +2.21%  "
$ws.Range("D34").Value = "'1.219"
$ws.Range("E34").Value = "'  +7.35%  "
$ws.Range("D35").Value = "'0.7585"
$ws.Range("E35").Value = "'  +2.65%  "
$ws.Range("E36").Value = "'  +0.65%  "
$ws.Range("D37").Value = "'0.01888"
$ws.Range("E37").Value = "'  +4.04%  "
$ws.Range("D39").Value = "'0.9222"
$ws.Range("D40").Value = "'2.081"
$ws.Range("E40").Value = "'  +1.38%  "
$ws.Range("D41").Value = "'107.03"
$ws.Range("E41").Value = "'  +0.27%  "
$ws.Range("D42").Value = "'0.4301"
$ws.Range("E42").Value = "'  +1.39%  "
$ws.Range("E43").Value = "'  +0.66%  "
$ws.Range("D44").Value = "'5.744"
$ws.Range("E44").Value = "'  -3.25%  "
$ws.Range("D45").Value = "'7.430"
$ws.Range("E45").Value = "'  +0.52%  "
$ws.Range("D46").Value = "'64.49"
$ws.Range("E46").Value = "'  +0.82%  "
$ws.Range("D47").Value = "'0.1276"
$ws.Range("E47").Value = "'  -2.66%  "
$ws.Range("D48").Value = "'1.486"
$ws.Range("E48").Value = "'  -4.96%  "
$ws.Range("D49").Value = "'8.957"
$ws.Range("E49").Value = "'  +2.08%  "
$ws.Range("D50").Value = "'33.87"
$ws.Range("D51").Value = "'0.05676"
$ws.Range("E51").Value = "'  +0.21%  "
